$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suivi")

# Insert a new column before CD, shifting CD->CE and CE->CF (matches xlShiftToRight)
$ws.Range("CD1").EntireColumn.Insert()

# Row 1 header: new timestamp cell, using same style as the rest of the header row
$ws.Range("CD1").Value = "2026-01-31 10:12:57"
$ws.Range("CD1").Style = $ws.Range("CC1").Style

# Body rows 2-80: carry forward last known price from column CC
$ws.Range("CD2:CD80").Value2 = $ws.Range("CC2:CC80").Value2

Write-Output "done"
